# Edit: add a new "Diretiva using: " paragraph (with an extra blank
# separator paragraph before it) right after the existing last
# paragraph ("MinhasClasses.console"), and relocate the hidden
# "_GoBack" bookmark into its own trailing empty paragraph.

$d = $word.ActiveDocument
$wNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# 1) The "_GoBack" bookmark currently lives inside the last paragraph
#    ("MinhasClasses.console"). Word keeps this bookmark hidden from
#    the normal Bookmarks enumeration, but it is still addressable by
#    name. Remove it from there - it will be re-created later in its
#    own paragraph at the very end of the document.
$goBackName = "_GoBack"
$hasGoBack = $false
foreach ($bm in @($goBackName)) {
    try {
        $existing = $d.Bookmarks($bm)
        $existing.Delete()
        $hasGoBack = $true
    } catch {
        $hasGoBack = $false
    }
}

# 2) Append three new paragraphs after the current last paragraph:
#      - an empty separator paragraph
#      - a paragraph that will hold "Diretiva using: "
#      - a paragraph that will hold the relocated bookmark
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

$blankPara = $d.Paragraphs($d.Paragraphs.Count)
$rBlank = $blankPara.Range
$rBlank.Collapse(0)
$rBlank.InsertParagraphAfter()

$diretivaParaTmp = $d.Paragraphs($d.Paragraphs.Count)
$rDiretivaTmp = $diretivaParaTmp.Range
$rDiretivaTmp.Collapse(0)
$rDiretivaTmp.InsertParagraphAfter()

# Re-fetch the three freshly created paragraphs by position.
$blankPara = $d.Paragraphs($d.Paragraphs.Count - 2)
$diretivaPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$bookmarkPara = $d.Paragraphs($d.Paragraphs.Count)

# 2a) Make sure the separator paragraph is a plain empty paragraph.
$rBlankFull = $blankPara.Range
$rBlankFull.InsertXML('<w:p ' + $wNS + '/>')

# 2b) Fill in the "Diretiva using: " paragraph, including the
#     spell-check proof markers ("using" is flagged as English inside
#     Portuguese text) around the word "using", matching the target
#     markup exactly.
$rDiretiva = $diretivaPara.Range
$diretivaXml = '<w:p ' + $wNS + '>' +
    '<w:r><w:t xml:space="preserve">Diretiva </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>using</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
    '</w:p>'
$rDiretiva.InsertXML($diretivaXml)

# 2c) Put the "_GoBack" bookmark into its own trailing paragraph.
$rBookmark = $bookmarkPara.Range
if ($hasGoBack) {
    $bookmarkXml = '<w:p ' + $wNS + '>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>'
    $rBookmark.InsertXML($bookmarkXml)
} else {
    $rBookmark.InsertXML('<w:p ' + $wNS + '/>')
}
